$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newAxioms = "SubClassOf: ns1:BFO_0000001 | SubClassOf: ies:Event | SubClassOf: ns1:BFO_0000132 min 0 | EquivalentTo: ies:Event"

# Row 2 <- values rotate in from old Row 4 (E,F,J); D gets new unified axioms; H unchanged
$ws.Range("D2").Value = $newAxioms
$ws.Range("E2").Value = "http://ies.data.gov.uk/ontology/ies4#Event"
$ws.Range("F2").Value = "Event"
$ws.Range("J2").Value = "An Event represents an activity or incident, involving one or more participating entities, that occurred/started at a specific point in time – e.g. a meeting, or a telephone call."

# Row 3 <- values rotate in from old Row 2 (E,F,H,J); D gets new unified axioms
$ws.Range("D3").Value = $newAxioms
$ws.Range("E3").Value = "http://ies.data.gov.uk/ontology/ies4#Entity"
$ws.Range("F3").Value = "Entity"
$ws.Range("H3").Value = "SubClassOf: ies:Element | SubClassOf: ies:isPartOf min 0"
$ws.Range("J3").Value = "An Entity typically represents a tangible thing like a Person, a Communications Device, or a Location."

# Row 4 <- values rotate in from old Row 3 (E,F,H,J); D gets new unified axioms
$ws.Range("D4").Value = $newAxioms
$ws.Range("E4").Value = "http://ies.data.gov.uk/ontology/ies4#State"
$ws.Range("F4").Value = "State"
$ws.Range("H4").Value = "SubClassOf: ies:isStateOf min 1"
$ws.Range("J4").Value = "A temporal state of an Element"
